$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[Sandro-Comandos Eletricos-1NB, Valmir-Caldeiraria-1NB, Anderson-Tornearia-1NB, Aselmo-Manut. Mot. End.-1NB]"
$ws.Range("C18").Value = "[Suzanny-Metalografia-1NB, Suzanny-Trat. Termicos-1NB, Gisele-E. D. N. D.-1NB, Weslei-Metrologia 1-1NB]"
$ws.Range("E18").Value = "[Andre B.-Elet. Dig. Bas.-1NB, Mayra-Tec. Mat. Não Metal.-1NB, Rachel-T.M. Metalicos-1NB, Victor-Ajustagem-1NB]"

$ws.Range("B19").Value = "[Sandro-Comandos Eletricos-1NB, Valmir-Caldeiraria-1NB, Anderson-Tornearia-1NB, Aselmo-Manut. Mot. End.-1NB]"
$ws.Range("C19").Value = "[Suzanny-Metalografia-1NB, Suzanny-Trat. Termicos-1NB, Gisele-E. D. N. D.-1NB, Weslei-Metrologia 1-1NB]"
$ws.Range("E19").Value = "[Andre B.-Elet. Dig. Bas.-1NB, Mayra-Tec. Mat. Não Metal.-1NB, Rachel-T.M. Metalicos-1NB, Victor-Ajustagem-1NB]"

$ws.Range("B20").Value = "[Sandro-Comandos Eletricos-1NB, Valmir-Caldeiraria-1NB, Anderson-Tornearia-1NB, Aselmo-Manut. Mot. End.-1NB]"
$ws.Range("C20").Value = "[Suzanny-Metalografia-1NB, Suzanny-Trat. Termicos-1NB, Gisele-E. D. N. D.-1NB, Weslei-Metrologia 1-1NB]"
$ws.Range("E20").Value = "[Andre B.-Elet. Dig. Bas.-1NB, Mayra-Tec. Mat. Não Metal.-1NB, Rachel-T.M. Metalicos-1NB, Victor-Ajustagem-1NB]"
$ws.Range("F20").Value = "Gilberto-M.T.R.M.-"

$ws.Range("B21").Value = "[Sandro-Comandos Eletricos-1NB, Valmir-Caldeiraria-1NB, Anderson-Tornearia-1NB, Aselmo-Manut. Mot. End.-1NB]"
$ws.Range("C21").Value = "[Suzanny-Metalografia-1NB, Suzanny-Trat. Termicos-1NB, Gisele-E. D. N. D.-1NB, Weslei-Metrologia 1-1NB]"
$ws.Range("E21").Value = "[Andre B.-Elet. Dig. Bas.-1NB, Mayra-Tec. Mat. Não Metal.-1NB, Rachel-T.M. Metalicos-1NB, Victor-Ajustagem-1NB]"
